# test02.xlsx ("Лист2"): mirror the Fibonacci sequence already in column A
# (rows 1-16) into column G, reproducing it cell-for-cell (seed values +
# relative "=prev2+prev1" formulas, auto-filled so Excel stores them as a
# shared formula group just like the original A column).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист2")

$ws.Range("G1").Value = 1
$ws.Range("G2").Value = 2
$ws.Range("G3").Formula = "=G2+G1"
$ws.Range("G4").Formula = "=G3+G2"
$ws.Range("G5:G16").Formula = "=G4+G3"

